$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = '"parent05'

$ws.Range("B5").Select()
